$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "In Translation"
#    This shared string is used by the Overview sheet (E2/F2 - one
#    status cell per locale) and by each locale sheet's "Status"
#    column (C2).
# ------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# ------------------------------------------------------------------
# 2) Narrow the "Status" columns to match the shorter text.
#    Target stored (OOXML) column width is 13.4101848602295
#    characters; ColumnWidth (the COM/UI measure) is the stored
#    width minus the standard 5px gridline padding expressed in
#    character units (0.8333333333333334 = 5/6 for the default
#    Calibri 11 font, MDW = 6px).
# ------------------------------------------------------------------
$targetColumnWidth = 13.4101848602295 - 0.8333333333333334

$wsOverview.Columns.Item(5).ColumnWidth = $targetColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $targetColumnWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $targetColumnWidth

$wsDeDe.Columns.Item(3).ColumnWidth = $targetColumnWidth
